$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume-1h (E) columns for changed rows
# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.136.40"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("E2").Style = "Normal"


# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.912.49"
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.02%  "
$ws.Range("E3").Style = "Normal"


# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("E4").Style = "Normal"


# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7395"
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.47%  "
$ws.Range("E5").Style = "Normal"


# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.85"
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("E6").Style = "Normal"


# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("E7").Style = "Normal"


# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.14%  "
$ws.Range("E8").Style = "Normal"


# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "26.92"
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.89%  "
$ws.Range("E9").Style = "Normal"


# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06985"
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.00%  "
$ws.Range("E10").Style = "Normal"


# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7800"
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.61%  "
$ws.Range("E11").Style = "Normal"


# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("E12").Style = "Normal"


# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.906.29"
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.43%  "
$ws.Range("E13").Style = "Normal"


# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.297"
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.61%  "
$ws.Range("E14").Style = "Normal"


# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.25"
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.19%  "
$ws.Range("E15").Style = "Normal"


# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.42"
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("E16").Style = "Normal"


# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.143.58"
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.66%  "
$ws.Range("E17").Style = "Normal"


# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.916"
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.98%  "
$ws.Range("E18").Style = "Normal"


# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.12"
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -4.20%  "
$ws.Range("E19").Style = "Normal"


# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007853"
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.76%  "
$ws.Range("E20").Style = "Normal"


# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("E21").Style = "Normal"


# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.152.49"
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.58%  "
$ws.Range("E22").Style = "Normal"


# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.33%  "
$ws.Range("E23").Style = "Normal"


# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.221"
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +8.58%  "
$ws.Range("E24").Style = "Normal"


# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.451"
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("E25").Style = "Normal"


# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.16"
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.38%  "
$ws.Range("E26").Style = "Normal"


# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.67%  "
$ws.Range("E27").Style = "Normal"


# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1287"
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.82%  "
$ws.Range("E28").Style = "Normal"


# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.070"
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -5.28%  "
$ws.Range("E29").Style = "Normal"


# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.354"
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.33%  "
$ws.Range("E30").Style = "Normal"


# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.29%  "
$ws.Range("E31").Style = "Normal"


# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.352"
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.75%  "
$ws.Range("E32").Style = "Normal"


# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.106"
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("E33").Style = "Normal"


# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05192"
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.22%  "
$ws.Range("E34").Style = "Normal"


# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.301"
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.69%  "
$ws.Range("E35").Style = "Normal"


# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7520"
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.52%  "
$ws.Range("E36").Style = "Normal"


# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.727"
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.67%  "
$ws.Range("E37").Style = "Normal"


# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01945"
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.43%  "
$ws.Range("E38").Style = "Normal"


# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.798"
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.22%  "
$ws.Range("E39").Style = "Normal"


# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.382"
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.63%  "
$ws.Range("E40").Style = "Normal"


# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4523"
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.83%  "
$ws.Range("E41").Style = "Normal"


# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "75.21"
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.56%  "
$ws.Range("E42").Style = "Normal"


# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.969"
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.68%  "
$ws.Range("E43").Style = "Normal"


# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.887"
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +5.67%  "
$ws.Range("E44").Style = "Normal"


# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("E45").Style = "Normal"


# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.8394"
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.80%  "
$ws.Range("E46").Style = "Normal"


# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "37.20"
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.53%  "
$ws.Range("E49").Style = "Normal"


# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.058.95"
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.33%  "
$ws.Range("E50").Style = "Normal"


# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1197"
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.99%  "
$ws.Range("E51").Style = "Normal"


# Row 47: now Quant (previously EnergySwap)
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "Quant"
$ws.Range("B47").Style = "Normal"

$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("C47").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.90"
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.48%  "
$ws.Range("E47").Style = "Normal"


# Row 48: now EnergySwap (previously Quant)
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("B48").Style = "Normal"

$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C48").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.961"
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.81%  "
$ws.Range("E48").Style = "Normal"


